$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the conversion text in A1 ---
$ws1 = $wb.Worksheets.Item("Hoja1")

$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 8.47 = 35005.08 pesos`n✅ 35005.08 pesos = 8.46 = 972.36 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

$ws1.Range("A1").Value = $newText

# --- Sheet "tasas": update the rate cells ---
$ws2 = $wb.Worksheets.Item("tasas")

$ws2.Range("N10").Value = 118.1
$ws2.Range("O10").Value = 4134.1
$ws2.Range("N12").Value = 4140
